$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1336.5652
$ws.Range("J17").Value = 1336.5652
$ws.Range("L17").Value = 4009.6956
$ws.Range("N17").Value = -4345.6956
$ws.Range("H33").Value = 559.64
$ws.Range("I33").Value = 371.94446
$ws.Range("J33").Value = 1042.2858
$ws.Range("K33").Value = 371.94446
$ws.Range("L33").Value = 1042.2858
$ws.Range("M33").Value = -142.94446
$ws.Range("N33").Value = -1500.2858
$ws.Range("H38").Value = 839.71875
$ws.Range("I38").Value = 202.27777
$ws.Range("J38").Value = 1659.2858
$ws.Range("K38").Value = 606.83331
$ws.Range("L38").Value = 4977.857400000001
$ws.Range("M38").Value = -234.83331
$ws.Range("N38").Value = -5721.857400000001
$ws.Range("H94").Value = 3250
$ws.Range("I94").Value = 3250
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 3250
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -2799
$ws.Range("N94").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1832.963
$ws.Range("I45").Value = 1091.25
$ws.Range("J45").Value = 7766.6665
$ws.Range("K45").Value = 1091.25
$ws.Range("L45").Value = 7766.6665
$ws.Range("M45").Value = -714.25
$ws.Range("N45").Value = -8520.666499999999
$ws.Range("H88").Value = 1595.3334
$ws.Range("I88").Value = 1595.3334
$ws.Range("K88").Value = 1595.3334
$ws.Range("M88").Value = -1189.3334
$ws.Range("H91").Value = 1595.3334
$ws.Range("I91").Value = 1595.3334
$ws.Range("K91").Value = 1595.3334
$ws.Range("M91").Value = -191.3334
$ws.Range("H122").Value = 4065.6365
$ws.Range("I122").Value = 2680.5
$ws.Range("J122").Value = 4857.143
$ws.Range("K122").Value = 8041.5
$ws.Range("L122").Value = 14571.429
$ws.Range("M122").Value = -5591.5
$ws.Range("N122").Value = -19471.429
$ws.Range("H132").Value = 26319880
$ws.Range("I132").Value = 40003780
$ws.Range("J132").Value = 4686.615
$ws.Range("K132").Value = 120011340
$ws.Range("L132").Value = 14059.845
$ws.Range("M132").Value = -120008810
$ws.Range("N132").Value = -19119.845

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3782.6667
$ws.Range("I20").Value = 1424
$ws.Range("J20").Value = 8500
$ws.Range("K20").Value = 1424
$ws.Range("L20").Value = 8500
$ws.Range("M20").Value = -1177
$ws.Range("N20").Value = -8994
$ws.Range("H134").Value = 2766
$ws.Range("I134").Value = 1645.5
$ws.Range("J134").Value = 5327.143
$ws.Range("K134").Value = 4936.5
$ws.Range("L134").Value = 15981.429
$ws.Range("M134").Value = -2401.5
$ws.Range("N134").Value = -21051.429

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H117").Value = 32000
$ws.Range("J117").Value = 32000
$ws.Range("L117").Value = 32000
$ws.Range("N117").Value = -41178
$ws.Range("H122").Value = 1853.4
$ws.Range("I122").Value = 2223.75
$ws.Range("J122").Value = 1679.1177
$ws.Range("K122").Value = 6671.25
$ws.Range("L122").Value = 5037.3531
$ws.Range("M122").Value = -4221.25
$ws.Range("N122").Value = -9937.3531

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H52").Value = 1084.4445
$ws.Range("J52").Value = 1084.4445
$ws.Range("L52").Value = 3253.3335
$ws.Range("N52").Value = -3785.3335
$ws.Range("H54").Value = 3600
$ws.Range("I54").Value = 2500
$ws.Range("J54").Value = 3820
$ws.Range("K54").Value = 7500
$ws.Range("L54").Value = 11460
$ws.Range("M54").Value = -6941
$ws.Range("N54").Value = -12578
$ws.Range("H60").Value = 366.66666
$ws.Range("I60").Value = 400
$ws.Range("J60").Value = 300
$ws.Range("K60").Value = 1200
$ws.Range("L60").Value = 900
$ws.Range("M60").Value = -949
$ws.Range("N60").Value = -1402
$ws.Range("H107").Value = 1397.2727
$ws.Range("J107").Value = 2330
$ws.Range("L107").Value = 6990
$ws.Range("N107").Value = -10830

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2234.4644
$ws.Range("I102").Value = 1640.579
$ws.Range("K102").Value = 1640.579
$ws.Range("M102").Value = -18.57899999999995
$ws.Range("H122").Value = 5244.3887
$ws.Range("I122").Value = 7766.5
$ws.Range("J122").Value = 3983.3333
$ws.Range("K122").Value = 23299.5
$ws.Range("L122").Value = 11949.9999
$ws.Range("M122").Value = -20849.5
$ws.Range("N122").Value = -16849.9999
$ws.Range("H126").Value = 2940
$ws.Range("I126").Value = 1486.6666
$ws.Range("J126").Value = 5120
$ws.Range("K126").Value = 4459.9998
$ws.Range("L126").Value = 15360
$ws.Range("M126").Value = -1989.9998
$ws.Range("N126").Value = -20300
$ws.Range("H132").Value = 3149.1143
$ws.Range("I132").Value = 2611.05
$ws.Range("J132").Value = 3866.5334
$ws.Range("K132").Value = 7833.150000000001
$ws.Range("L132").Value = 11599.6002
$ws.Range("M132").Value = -5303.150000000001
$ws.Range("N132").Value = -16659.6002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 429.4
$ws.Range("I16").Value = 412.375
$ws.Range("K16").Value = 412.375
$ws.Range("M16").Value = -242.375
$ws.Range("H55").Value = 766.2105
$ws.Range("I55").Value = 108
$ws.Range("J55").Value = 1244.909
$ws.Range("K55").Value = 108
$ws.Range("L55").Value = 1244.909
$ws.Range("M55").Value = 65
$ws.Range("N55").Value = -1590.909
$ws.Range("H93").Value = 4360.5
$ws.Range("I93").Value = 3472.6365
$ws.Range("J93").Value = 5755.7144
$ws.Range("K93").Value = 3472.6365
$ws.Range("L93").Value = 5755.7144
$ws.Range("M93").Value = -2224.6365
$ws.Range("N93").Value = -8251.714400000001
$ws.Range("H111").Value = 39455.332
$ws.Range("J111").Value = 39455.332
$ws.Range("L111").Value = 39455.332
$ws.Range("N111").Value = -47635.332
$ws.Range("H132").Value = 2154.34
$ws.Range("I132").Value = 1139.9697
$ws.Range("J132").Value = 4123.4116
$ws.Range("K132").Value = 3419.9091
$ws.Range("L132").Value = 12370.2348
$ws.Range("M132").Value = -889.9091000000003
$ws.Range("N132").Value = -17430.2348

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H98").Value = 22999.5
$ws.Range("J98").Value = 22999.5
$ws.Range("L98").Value = 22999.5
$ws.Range("N98").Value = -28989.5
$ws.Range("H122").Value = 401875.03
$ws.Range("I122").Value = 557239.6
$ws.Range("J122").Value = 2366.1428
$ws.Range("K122").Value = 1671718.8
$ws.Range("L122").Value = 7098.428400000001
$ws.Range("M122").Value = -1669268.8
$ws.Range("N122").Value = -11998.4284
$ws.Range("H132").Value = 19297.383
$ws.Range("I132").Value = 5127.636
$ws.Range("J132").Value = 45275.25
$ws.Range("K132").Value = 15382.908
$ws.Range("L132").Value = 135825.75
$ws.Range("M132").Value = -12852.908
$ws.Range("N132").Value = -140885.75
